# "Add files via upload" — re-saved workbook: the OLED wiring block (rows
# 17:20) on the "VERSION 2!!!" sheet is removed, leaving only the two
# formatted-but-empty rows behind, and a couple of cosmetic/formatting
# touch-ups carried along with the re-save.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("VERSION 2!!!")

# Remove the two fully-populated OLED rows (19:20) outright...
$ws2.Rows("19:20").Delete()

# ...and blank out the remaining two rows' contents (B:E), leaving the
# already-formatted D17/D18 cells behind as empty styled cells.
$ws2.Range("B17:E18").ClearContents()

# Reflect where the selection ended up after the edit.
$ws2.Activate()
$ws2.Range("F19").Select()

# Row-height tweaks on Version1 that came along with the re-save.
$ws1 = $wb.Worksheets.Item("Version1")
$ws1.Rows.Item(8).RowHeight = 12.9
$ws1.Rows.Item(10).RowHeight = 12.9

# Default workbook font switched from the Korean "맑은 고딕" to "Calibri"
# for the Normal / Good / Bad cell styles.
$wb.Styles.Item("Normal").Font.Name = "Calibri"
$wb.Styles.Item("Good").Font.Name = "Calibri"
$wb.Styles.Item("Bad").Font.Name = "Calibri"
